# Update "想去人数" (interest count) figures in the 展览 (Exhibitions) sheet
# and the combined 全部类型 (All types) sheet, per the upstream data refresh.
#
#   F3: 86   -> 88
#   F4: 5    -> 11
#   F5: 2269 -> 2287
#   F6: 201  -> 204
#   F7/F9: 370 -> 371   (row 7 on 展览, row 9 on 全部类型)

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3" = 88
    "F4" = 11
    "F5" = 2287
    "F6" = 204
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}

# Last row differs between the two sheets: row 7 on 展览, row 9 on 全部类型.
$wb.Worksheets.Item("展览").Range("F7").Value = 371
$wb.Worksheets.Item("全部类型").Range("F9").Value = 371
